# Weekly crime data update for 123rd Precinct CompStat report
# - Updates masthead: Police Commissioner name, Volume/Number, and report date range
# - Updates the Crime Complaints grid (rows 16-21, 24-26, 28) with newly collected figures

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Masthead text updates -------------------------------------------------
$ws.Range("M6").Value = "Jessica S. Tisch"
$ws.Range("A8").Value = "Volume 31   Number  48"
$ws.Range("C9").Value = "Report Covering the Week  11/25/2024  Through  12/1/2024"

# --- Crime Complaints grid updates (rows 16-21, 24-26, 28) -----------------
# Some cells switch between the "not applicable" placeholder text (styled
# General, showing "0" or "***.*") and real numeric/percentage values, so a
# handful of cells first have their number format copied from a same-column
# reference cell (row 14, which keeps its original layout) before the new
# value is written, to keep the cell style consistent with the new content.
$ws.Range("C16").Value = 2
$ws.Range("I14").Copy()
$ws.Range("D16").PasteSpecial(-4122)
$ws.Range("D16").Value = 2
$ws.Range("K14").Copy()
$ws.Range("E16").PasteSpecial(-4122)
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 5
$ws.Range("G16").Value = 3
$ws.Range("H16").Value = 66.666666666666
$ws.Range("I16").Value = 25
$ws.Range("J16").Value = 20
$ws.Range("K16").Value = 25
$ws.Range("L16").Value = 38.888888888888
$ws.Range("M16").Value = 13.636363636363
$ws.Range("N16").Value = -59.677419354838
$ws.Range("I14").Copy()
$ws.Range("C17").PasteSpecial(-4122)
$ws.Range("C17").Value = 2
$ws.Range("D17").Value = 2
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 8
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 86
$ws.Range("J17").Value = 84
$ws.Range("K17").Value = 2.380952380952
$ws.Range("L17").Value = 68.627450980392
$ws.Range("M17").Value = 95.454545454545
$ws.Range("N17").Value = -15.686274509803
$ws.Range("C18").Value = 1
$ws.Range("I14").Copy()
$ws.Range("D18").PasteSpecial(-4122)
$ws.Range("D18").Value = 2
$ws.Range("K14").Copy()
$ws.Range("E18").PasteSpecial(-4122)
$ws.Range("E18").Value = -50
$ws.Range("F18").Value = 7
$ws.Range("G18").Value = 4
$ws.Range("H18").Value = 75
$ws.Range("I18").Value = 41
$ws.Range("J18").Value = 51
$ws.Range("K18").Value = -19.607843137254
$ws.Range("L18").Value = 2.5
$ws.Range("M18").Value = -58.163265306122
$ws.Range("N18").Value = -86.858974358974
$ws.Range("C19").Value = 2
$ws.Range("E19").Value = -50
$ws.Range("F19").Value = 9
$ws.Range("G19").Value = 23
$ws.Range("H19").Value = -60.869565217391
$ws.Range("I19").Value = 242
$ws.Range("J19").Value = 268
$ws.Range("K19").Value = -9.701492537313
$ws.Range("L19").Value = -2.419354838709
$ws.Range("M19").Value = 77.941176470588
$ws.Range("N19").Value = 8.035714285714
$ws.Range("I14").Copy()
$ws.Range("D20").PasteSpecial(-4122)
$ws.Range("D20").Value = 1
$ws.Range("K14").Copy()
$ws.Range("E20").PasteSpecial(-4122)
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 3
$ws.Range("G20").Value = 3
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 37
$ws.Range("J20").Value = 73
$ws.Range("K20").Value = -49.315068493150
$ws.Range("L20").Value = -67.543859649122
$ws.Range("M20").Value = 2.777777777777
$ws.Range("N20").Value = -94.359756097561
$ws.Range("D21").Value = 11
$ws.Range("E21").Value = -27.272727272727
$ws.Range("F21").Value = 32
$ws.Range("G21").Value = 41
$ws.Range("H21").Value = -21.951219512195
$ws.Range("I21").Value = 439
$ws.Range("J21").Value = 498
$ws.Range("K21").Value = -11.847389558232
$ws.Range("L21").Value = -7.578947368421
$ws.Range("M21").Value = 29.117647058823
$ws.Range("N21").Value = -67.744305657604
$ws.Range("C24").Value = 9
$ws.Range("D24").Value = 16
$ws.Range("E24").Value = -43.75
$ws.Range("F24").Value = 47
$ws.Range("G24").Value = 50
$ws.Range("H24").Value = -6
$ws.Range("I24").Value = 351
$ws.Range("J24").Value = 460
$ws.Range("K24").Value = -23.695652173913
$ws.Range("L24").Value = -22.857142857142
$ws.Range("M24").Value = -30.079681274900
$ws.Range("C25").Value = 5
$ws.Range("D25").Value = 4
$ws.Range("E25").Value = 25
$ws.Range("F25").Value = 28
$ws.Range("G25").Value = 12
$ws.Range("H25").Value = 133.333333333333
$ws.Range("I25").Value = 163
$ws.Range("J25").Value = 191
$ws.Range("K25").Value = -14.659685863874
$ws.Range("L25").Value = 45.535714285714
$ws.Range("C26").Value = 2
$ws.Range("E26").Value = -60
$ws.Range("F26").Value = 15
$ws.Range("G26").Value = 14
$ws.Range("H26").Value = 7.142857142857
$ws.Range("I26").Value = 152
$ws.Range("J26").Value = 180
$ws.Range("K26").Value = -15.555555555555
$ws.Range("L26").Value = -11.111111111111
$ws.Range("M26").Value = -23.232323232323
$ws.Range("I14").Copy()
$ws.Range("D28").PasteSpecial(-4122)
$ws.Range("D28").Value = 1
$ws.Range("K14").Copy()
$ws.Range("E28").PasteSpecial(-4122)
$ws.Range("E28").Value = -100
$ws.Range("C14").Copy()
$ws.Range("F28").PasteSpecial(-4122)
$ws.Range("F28").NumberFormat = "@"
$ws.Range("F28").Value = "0"
$ws.Range("I14").Copy()
$ws.Range("G28").PasteSpecial(-4122)
$ws.Range("G28").Value = 1
$ws.Range("K14").Copy()
$ws.Range("H28").PasteSpecial(-4122)
$ws.Range("H28").Value = -100
$ws.Range("J28").Value = 17
$ws.Range("K28").Value = 17.647058823529

$excel.CutCopyMode = $false
